# Update answer values in the division practice table.
# Each pair below is (old expression => new expression); replacements are
# unique substrings so Find/Replace targets exactly one cell each.
$d = $word.ActiveDocument

$replacements = @(
    @("337÷6=56, 1", "435÷4=108, 3"),
    @("972÷3=324, 0", "384÷6=64, 0"),
    @("493÷5=98, 3", "764÷2=382, 0"),
    @("919÷6=153, 1", "282÷2=141, 0"),
    @("282÷3=94, 0", "318÷4=79, 2"),
    @("570÷5=114, 0", "658÷5=131, 3"),
    @("338÷3=112, 2", "491÷6=81, 5"),
    @("672÷9=74, 6", "237÷7=33, 6"),
    @("134÷3=44, 2", "589÷6=98, 1"),
    @("868÷5=173, 3", "130÷5=26, 0"),
    @("770÷9=85, 5", "811÷7=115, 6"),
    @("121÷9=13, 4", "249÷5=49, 4"),
    @("396÷2=198, 0", "406÷3=135, 1"),
    @("992÷8=124, 0", "530÷7=75, 5"),
    @("394÷7=56, 2", "732÷3=244, 0"),
    @("969÷8=121, 1", "819÷6=136, 3"),
    @("882÷8=110, 2", "574÷8=71, 6"),
    @("469÷5=93, 4", "316÷3=105, 1"),
    @("135÷4=33, 3", "968÷3=322, 2"),
    @("804÷8=100, 4", "498÷6=83, 0"),
    @("910÷8=113, 6", "462÷8=57, 6"),
    @("327÷9=36, 3", "933÷9=103, 6"),
    @("256÷7=36, 4", "654÷4=163, 2"),
    @("847÷3=282, 1", "349÷7=49, 6"),
    @("468÷8=58, 4", "604÷5=120, 4"),
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Could not find text to replace: $oldText"
    }
}
